$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.48080775141716
$ws.Range("B1").Value = 1.714520215988159
$ws.Range("C1").Value = 2.278101921081543
$ws.Range("D1").Value = 1.959650278091431
$ws.Range("E1").Value = 0.9656597375869751
